$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("AH4").Value = 51
$ws.Range("AN4").Value = 11
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 2.55
$ws.Range("Y4").Value = 1.67
$ws.Range("Z4").Value = 2.1
# Row 5
$ws.Range("AD5").Value = 13
$ws.Range("AF5").Value = 34
$ws.Range("AM5").Value = 6
$ws.Range("G5").Value = 3.25
$ws.Range("I5").Value = 2.45
$ws.Range("J5").Value = 4
# Row 7
$ws.Range("AJ7").Value = 8
$ws.Range("K7").Value = 2.3
# Row 8
$ws.Range("AA8").Value = 1.67
$ws.Range("AB8").Value = 2.1
$ws.Range("AC8").Value = 8.5
$ws.Range("AD8").Value = 11
$ws.Range("AF8").Value = 19
$ws.Range("AG8").Value = 17
$ws.Range("AI8").Value = 11
$ws.Range("AM8").Value = 11
$ws.Range("AN8").Value = 17
$ws.Range("AO8").Value = 12
$ws.Range("AP8").Value = 34
$ws.Range("G8").Value = 2.05
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 3.3
$ws.Range("J8").Value = 2.75
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 3.75
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("S8").Value = 1.85
$ws.Range("T8").Value = 2
# Row 9
$ws.Range("AF9").Value = 12
$ws.Range("AI9").Value = 15
$ws.Range("AM9").Value = 19
$ws.Range("AP9").Value = 51
$ws.Range("G9").Value = 1.53
$ws.Range("H9").Value = 4.2
$ws.Range("I9").Value = 5.75
$ws.Range("J9").Value = 2.05
$ws.Range("S9").Value = 1.57
$ws.Range("T9").Value = 2.35
$ws.Range("U9").Value = 1.9
$ws.Range("V9").Value = 1.9
$ws.Range("W9").Value = 2.38
$ws.Range("X9").Value = 1.53
# Row 10
$ws.Range("AA10").Value = 1.8
$ws.Range("AB10").Value = 1.95
$ws.Range("AG10").Value = 21
$ws.Range("AI10").Value = 9
$ws.Range("AM10").Value = 8.5
$ws.Range("AN10").Value = 13
$ws.Range("AS10").Value = 251
$ws.Range("G10").Value = 2.7
$ws.Range("H10").Value = 3.25
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 3.25
$ws.Range("S10").Value = 2.08
$ws.Range("T10").Value = 1.73
$ws.Range("W10").Value = 3.75
$ws.Range("X10").Value = 1.25
# Row 11
$ws.Range("AA11").Value = 1.75
$ws.Range("AB11").Value = 2
$ws.Range("AC11").Value = 8.5
$ws.Range("AH11").Value = 21
$ws.Range("AI11").Value = 15
$ws.Range("AJ11").Value = 8.5
$ws.Range("AL11").Value = 41
$ws.Range("G11").Value = 1.53
$ws.Range("H11").Value = 4.5
$ws.Range("I11").Value = 5.5
$ws.Range("J11").Value = 2.05
$ws.Range("K11").Value = 2.4
$ws.Range("M11").Value = 1.03
$ws.Range("N11").Value = 15
$ws.Range("O11").Value = 1.18
$ws.Range("P11").Value = 4.5
$ws.Range("S11").Value = 1.62
$ws.Range("T11").Value = 2.25
$ws.Range("W11").Value = 2.5
$ws.Range("X11").Value = 1.5
$ws.Range("Y11").Value = 1.3
$ws.Range("Z11").Value = 3.4
# Row 12
$ws.Range("AA12").Value = 1.95
$ws.Range("AB12").Value = 1.8
$ws.Range("AC12").Value = 9
$ws.Range("AI12").Value = 8
$ws.Range("AR12").Value = 34
$ws.Range("AS12").Value = 351
$ws.Range("G12").Value = 3.6
$ws.Range("H12").Value = 3.1
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 9
$ws.Range("O12").Value = 1.36
$ws.Range("P12").Value = 3
$ws.Range("S12").Value = 2.2
$ws.Range("T12").Value = 1.65
$ws.Range("Y12").Value = 1.5
$ws.Range("Z12").Value = 2.5
# Row 13
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 10
$ws.Range("T13").Value = 1.72
# Row 14
$ws.Range("S14").Value = 1.9
$ws.Range("T14").Value = 1.95
# Row 15
$ws.Range("N15").Value = 8
$ws.Range("S15").Value = 2.3
$ws.Range("T15").Value = 1.6
$ws.Range("W15").Value = 4.33
$ws.Range("X15").Value = 1.2
# Row 17
$ws.Range("AA17").Value = 1.37
$ws.Range("AD17").Value = 18
$ws.Range("AI17").Value = 9.75
$ws.Range("AJ17").Value = 8
$ws.Range("AP17").Value = 32
$ws.Range("G17").Value = 2.52
$ws.Range("H17").Value = 3.6
$ws.Range("J17").Value = 2.92
$ws.Range("K17").Value = 2.37
$ws.Range("N17").Value = 9.75
$ws.Range("T17").Value = 2.67
$ws.Range("Z17").Value = 3.6
# Row 18
$ws.Range("AB18").Value = 2.32
$ws.Range("AC18").Value = 17
$ws.Range("AD18").Value = 29
$ws.Range("AF18").Value = 75
$ws.Range("AI18").Value = 9.25
$ws.Range("AL18").Value = 40
$ws.Range("AO18").Value = 8.25
$ws.Range("AP18").Value = 14
$ws.Range("G18").Value = 4.45
$ws.Range("H18").Value = 3.95
$ws.Range("K18").Value = 2.4
$ws.Range("N18").Value = 9.25
$ws.Range("P18").Value = 4.5
$ws.Range("Y18").Value = 1.28
$ws.Range("Z18").Value = 3.35
# Row 20
$ws.Range("AA20").Value = 1.75
$ws.Range("AB20").Value = 1.95
$ws.Range("AC20").Value = 17
$ws.Range("AD20").Value = 35
$ws.Range("AF20").Value = 110
$ws.Range("AH20").Value = 45
$ws.Range("AK20").Value = 15.5
$ws.Range("AN20").Value = 7.6
$ws.Range("AS20").Value = 450
$ws.Range("G20").Value = 5.5
$ws.Range("H20").Value = 4.05
$ws.Range("J20").Value = 5.3
$ws.Range("K20").Value = 2.3
$ws.Range("P20").Value = 3.85
$ws.Range("W20").Value = 2.57
$ws.Range("Y20").Value = 1.34
$ws.Range("Z20").Value = 3
